$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HW-Net")
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("Q3").Select()
